# Refresh the crypto price/volume snapshot (GitHub Actions data pull).
# Also corrects the ranking order for Avalanche/BitDAO (rows 22-23 swap places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'22.383.93"
$ws.Range('E2').Value = '  -4.14%  '

$ws.Range('D3').Value = "'1.571.83"
$ws.Range('E3').Value = '  -3.70%  '

$ws.Range('D4').Value = "'1.000"
$ws.Range('E4').Value = '  -0.24%  '

$ws.Range('D5').Value = "'1.000"
$ws.Range('E5').Value = '  -0.21%  '

$ws.Range('D6').Value = "'290.02"
$ws.Range('E6').Value = '  -3.17%  '

$ws.Range('D7').Value = "'0.3692"
$ws.Range('E7').Value = '  -2.19%  '

$ws.Range('D8').Value = "'49.42"
$ws.Range('E8').Value = '  -1.58%  '

$ws.Range('D9').Value = "'0.3384"
$ws.Range('E9').Value = '  -3.70%  '

$ws.Range('D10').Value = "'1.166"
$ws.Range('E10').Value = '  -3.10%  '

$ws.Range('D11').Value = "'0.07607"
$ws.Range('E11').Value = '  -5.37%  '

$ws.Range('E12').Value = '  -0.22%  '

$ws.Range('D13').Value = "'21.28"
$ws.Range('E13').Value = '  -2.66%  '

$ws.Range('D14').Value = "'6.054"
$ws.Range('E14').Value = '  -3.95%  '

$ws.Range('D15').Value = "'6.904"
$ws.Range('E15').Value = '  -4.55%  '

$ws.Range('D16').Value = "'1.572.69"
$ws.Range('E16').Value = '  -3.99%  '

$ws.Range('D17').Value = "'0.00001132"
$ws.Range('E17').Value = '  -5.54%  '

$ws.Range('D18').Value = "'89.27"
$ws.Range('E18').Value = '  -6.71%  '

$ws.Range('D19').Value = "'0.06759"
$ws.Range('E19').Value = '  -2.88%  '

$ws.Range('E20').Value = '  -0.09%  '

$ws.Range('D21').Value = "'6.248"
$ws.Range('E21').Value = '  -6.61%  '

$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').Value = "'16.55"
$ws.Range('E22').Value = '  -4.07%  '

$ws.Range('B23').Value = 'BitDAO'
$ws.Range('C23').Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range('D23').Value = "'0.5307"
$ws.Range('E23').Value = '  -7.70%  '

$ws.Range('D24').Value = "'11.96"
$ws.Range('E24').Value = '  -2.45%  '

$ws.Range('D25').Value = "'22.387.62"
$ws.Range('E25').Value = '  -4.18%  '

$ws.Range('D26').Value = "'2.388"
$ws.Range('E26').Value = '  -3.23%  '

$ws.Range('D27').Value = "'3.007"
$ws.Range('E27').Value = '  +4.86%  '

$ws.Range('D28').Value = "'19.95"
$ws.Range('E28').Value = '  -3.86%  '

$ws.Range('D29').Value = "'145.54"
$ws.Range('E29').Value = '  -4.21%  '

$ws.Range('D30').Value = "'4.973"
$ws.Range('E30').Value = '  -3.93%  '

$ws.Range('D31').Value = "'125.60"
$ws.Range('E31').Value = '  -4.88%  '

$ws.Range('D32').Value = "'1.742.31"
$ws.Range('E32').Value = '  -4.29%  '

$ws.Range('D33').Value = "'1.048"
$ws.Range('E33').Value = '  +8.20%  '

$ws.Range('D34').Value = "'6.272"
$ws.Range('E34').Value = '  -7.59%  '

$ws.Range('D35').Value = "'1.994"
$ws.Range('E35').Value = '  -6.29%  '

$ws.Range('D36').Value = "'10.34"
$ws.Range('E36').Value = '  -8.45%  '

$ws.Range('D37').Value = "'0.08508"
$ws.Range('E37').Value = '  -2.36%  '

$ws.Range('D38').Value = "'0.02533"
$ws.Range('E38').Value = '  -5.99%  '

$ws.Range('D39').Value = "'0.2328"
$ws.Range('E39').Value = '  -3.63%  '

$ws.Range('D40').Value = "'5.564"
$ws.Range('E40').Value = '  -4.81%  '

$ws.Range('D41').Value = "'0.06516"
$ws.Range('E41').Value = '  -3.84%  '

$ws.Range('D42').Value = "'11.80"
$ws.Range('E42').Value = '  -8.30%  '

$ws.Range('D43').Value = "'1.244"
$ws.Range('E43').Value = '  -4.12%  '

$ws.Range('D44').Value = "'0.6371"
$ws.Range('E44').Value = '  -6.36%  '

$ws.Range('D45').Value = "'14.37"
$ws.Range('E45').Value = '  -8.02%  '

$ws.Range('D46').Value = "'0.9999"
$ws.Range('E46').Value = '  -0.08%  '

$ws.Range('D47').Value = "'0.5991"
$ws.Range('E47').Value = '  -4.73%  '

$ws.Range('D48').Value = "'3.757"
$ws.Range('E48').Value = '  -3.56%  '

$ws.Range('D49').Value = "'2.128"
$ws.Range('E49').Value = '  -4.56%  '

$ws.Range('D50').Value = "'1.253"
$ws.Range('E50').Value = '  +4.15%  '

$ws.Range('D51').Value = "'123.44"
$ws.Range('E51').Value = '  -2.50%  '
